$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; this shifts the existing rows 23..30
# down to 24..31 and the worksheet's used range grows from A1:R30 to A1:R31.
$ws.Rows.Item(23).Insert()

# The new row at position 23 keeps the same constant columns as its
# neighbours (Mercado, Region, Codreg, Categoria ID, etc.) and gets the
# week's new price data.
$ws.Cells.Item(23, 1).Value = 7
$ws.Cells.Item(23, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(23, 3).Value = "Ñuble"
$ws.Cells.Item(23, 4).Value = 44524
$ws.Cells.Item(23, 4).NumberFormat = $ws.Cells.Item(24, 4).NumberFormat
$ws.Cells.Item(23, 5).Value = 16
$ws.Cells.Item(23, 6).Value = 100112026
$ws.Cells.Item(23, 7).Value = "Haba"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 6000
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 6500
$ws.Cells.Item(23, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(23, 16).Value = 260
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
